$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1895.125
$ws.Range("I40").Value = 1860.1666
$ws.Range("K40").Value = 1860.1666
$ws.Range("M40").Value = -1685.1666
$ws.Range("H64").Value = 3641.25
$ws.Range("I64").Value = 3187.3333
$ws.Range("J64").Value = 5003
$ws.Range("K64").Value = 3187.3333
$ws.Range("L64").Value = 5003
$ws.Range("M64").Value = -2939.3333
$ws.Range("N64").Value = -5499
$ws.Range("H67").Value = 3641.25
$ws.Range("I67").Value = 3187.3333
$ws.Range("J67").Value = 5003
$ws.Range("K67").Value = 3187.3333
$ws.Range("L67").Value = 5003
$ws.Range("M67").Value = -2329.3333
$ws.Range("N67").Value = -6719
$ws.Range("H137").Value = 985.381
$ws.Range("I137").Value = 882.5294
$ws.Range("K137").Value = 2647.5882
$ws.Range("M137").Value = -97.58820000000014
$ws.Range("H140").Value = 64607.5
$ws.Range("J140").Value = 64607.5
$ws.Range("L140").Value = 64607.5
$ws.Range("N140").Value = -74967.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1947.4736
$ws.Range("I61").Value = 1152.1818
$ws.Range("J61").Value = 3041
$ws.Range("K61").Value = 1152.1818
$ws.Range("L61").Value = 3041
$ws.Range("M61").Value = -940.1818000000001
$ws.Range("N61").Value = -3465
$ws.Range("H74").Value = 1134.0938
$ws.Range("I74").Value = 914.45
$ws.Range("K74").Value = 914.45
$ws.Range("M74").Value = -40.45000000000005
$ws.Range("H77").Value = 1134.0938
$ws.Range("I77").Value = 914.45
$ws.Range("K77").Value = 4572.25
$ws.Range("M77").Value = -204.25
$ws.Range("H132").Value = 15641901
$ws.Range("I132").Value = 21277640
$ws.Range("J132").Value = 60740.234
$ws.Range("K132").Value = 63832920
$ws.Range("L132").Value = 182220.702
$ws.Range("M132").Value = -63830390
$ws.Range("N132").Value = -187280.702
$ws.Range("H136").Value = 1947.4736
$ws.Range("I136").Value = 1152.1818
$ws.Range("J136").Value = 3041
$ws.Range("K136").Value = 3456.5454
$ws.Range("L136").Value = 9123
$ws.Range("M136").Value = -906.5454
$ws.Range("N136").Value = -14223

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3940.2
$ws.Range("I105").Value = 4422.5
$ws.Range("J105").Value = 2011
$ws.Range("K105").Value = 4422.5
$ws.Range("L105").Value = 2011
$ws.Range("M105").Value = -2675.5
$ws.Range("N105").Value = -5505
$ws.Range("H134").Value = 2093.5715
$ws.Range("I134").Value = 2030.5
$ws.Range("J134").Value = 2766.3333
$ws.Range("K134").Value = 6091.5
$ws.Range("L134").Value = 8298.999899999999
$ws.Range("M134").Value = -3556.5
$ws.Range("N134").Value = -13368.9999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4106.8057
$ws.Range("I31").Value = 5712.773
$ws.Range("J31").Value = 1583.1428
$ws.Range("K31").Value = 5712.773
$ws.Range("L31").Value = 1583.1428
$ws.Range("M31").Value = -5417.773
$ws.Range("N31").Value = -2173.1428
$ws.Range("H34").Value = 4106.8057
$ws.Range("I34").Value = 5712.773
$ws.Range("J34").Value = 1583.1428
$ws.Range("K34").Value = 5712.773
$ws.Range("L34").Value = 1583.1428
$ws.Range("M34").Value = -5510.773
$ws.Range("N34").Value = -1987.1428
$ws.Range("H58").Value = 3146.9473
$ws.Range("I58").Value = 1298.6666
$ws.Range("K58").Value = 1298.6666
$ws.Range("M58").Value = -1095.6666
$ws.Range("H132").Value = 48606.41
$ws.Range("I132").Value = 2289.4666
$ws.Range("K132").Value = 6868.399800000001
$ws.Range("M132").Value = -4338.399800000001
$ws.Range("H134").Value = 3439.5454
$ws.Range("I134").Value = 2119.2856
$ws.Range("J134").Value = 5750
$ws.Range("K134").Value = 6357.8568
$ws.Range("L134").Value = 17250
$ws.Range("M134").Value = -3822.8568
$ws.Range("N134").Value = -22320
$ws.Range("H136").Value = 3146.9473
$ws.Range("I136").Value = 1298.6666
$ws.Range("K136").Value = 3895.9998
$ws.Range("M136").Value = -1345.9998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 1227
$ws.Range("I40").Value = 104.117645
$ws.Range("J40").Value = 5999.25
$ws.Range("K40").Value = 416.47058
$ws.Range("L40").Value = 23997
$ws.Range("M40").Value = -347.47058
$ws.Range("N40").Value = -24135
$ws.Range("H68").Value = 1400
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 1560
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 4680
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -6302
$ws.Range("H71").Value = 1400
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 1560
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 14040
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -22152
$ws.Range("H75").Value = 750
$ws.Range("I75").Value = 700
$ws.Range("K75").Value = 2100
$ws.Range("M75").Value = -1102
$ws.Range("N75").Value = -4396
$ws.Range("H76").Value = 4920
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 9000
$ws.Range("L76").Value = 15000
$ws.Range("M76").Value = -8617
$ws.Range("N76").Value = -15766
$ws.Range("H78").Value = 750
$ws.Range("I78").Value = 700
$ws.Range("K78").Value = 6300
$ws.Range("M78").Value = -1308
$ws.Range("N78").Value = -17184
$ws.Range("H79").Value = 4920
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 9000
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = -7674
$ws.Range("N79").Value = -17652
$ws.Range("H80").Value = 2375
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 7500
$ws.Range("M80").Value = -5064
$ws.Range("N80").Value = -9372
$ws.Range("H81").Value = 4172.222
$ws.Range("I81").Value = 190
$ws.Range("J81").Value = 4490.8
$ws.Range("K81").Value = 570
$ws.Range("L81").Value = 13472.4
$ws.Range("M81").Value = 553
$ws.Range("N81").Value = -15718.4
$ws.Range("H83").Value = 2375
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 22500
$ws.Range("M83").Value = -13320
$ws.Range("N83").Value = -31860
$ws.Range("H84").Value = 4172.222
$ws.Range("I84").Value = 190
$ws.Range("J84").Value = 4490.8
$ws.Range("K84").Value = 1710
$ws.Range("L84").Value = 40417.2
$ws.Range("M84").Value = 3906
$ws.Range("N84").Value = -51649.2
$ws.Range("H107").Value = 328.8846
$ws.Range("I107").Value = 510
$ws.Range("J107").Value = 305.26086
$ws.Range("K107").Value = 1530
$ws.Range("L107").Value = 915.7825799999999
$ws.Range("M107").Value = 390
$ws.Range("N107").Value = -4755.78258
$ws.Range("H120").Value = 14009.667
$ws.Range("I120").Value = 11014.5
$ws.Range("J120").Value = 20000
$ws.Range("K120").Value = 33043.5
$ws.Range("L120").Value = 60000
$ws.Range("M120").Value = -28205.5
$ws.Range("N120").Value = -69676

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 37916.855
$ws.Range("I132").Value = 1504.2858
$ws.Range("K132").Value = 4512.857400000001
$ws.Range("M132").Value = -1982.857400000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 12389.333
$ws.Range("I136").Value = 12286.546
$ws.Range("K136").Value = 36859.638
$ws.Range("M136").Value = -34309.638
$ws.Range("H140").Value = 142912340
$ws.Range("J140").Value = 61865.8
$ws.Range("L140").Value = 61865.8
$ws.Range("N140").Value = -72225.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 84900390
$ws.Range("I132").Value = 161429840
$ws.Range("K132").Value = 484289520
$ws.Range("M132").Value = -484286990
$ws.Range("H136").Value = 23118.133
$ws.Range("I136").Value = 42429.625
$ws.Range("J136").Value = 1047.8572
$ws.Range("K136").Value = 127288.875
$ws.Range("L136").Value = 3143.5716
$ws.Range("M136").Value = -124738.875
$ws.Range("N136").Value = -8243.571599999999

Write-Host "Applied all profit/price updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."